# Add a new scenario row ("ieaghg-reference") to the scenario library sheet,
# matching the IEAGHG reference-scenario update described in the commit
# message ("updates to IEAGHG scenario / rechecking of emission and energy
# data based on source").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Range("A$row").Value = "ieaghg-reference"
$ws.Range("B$row").Value = "coal"
$ws.Range("C$row").Value = "air"
$ws.Range("D$row").Value = "BF production"
$ws.Range("E$row").Value = "electricity partially co-generated from fluegases"

# Technology Notes column wraps text (same style as the other rows) and the
# row height matches the other note rows (16pt).
$ws.Range("D$row").WrapText = $true
$ws.Rows.Item($row).RowHeight = 16
